# The commit republishes the SectorGroup codelist with columns E
# ("codeforiati:category-name") and F ("codeforiati:group-code") swapped -
# the header labels trade places and, for every data row, the category-name
# text and the group-code text trade places too (columns D and G are left
# untouched). Net effect: column E now holds the group-code, column F now
# holds the category-name, for the header row and every data row alike.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count()

# Use Copy/PasteSpecial (through an unused scratch column) rather than a
# plain Value assignment: round-tripping through .Value would make Excel
# re-infer the type of purely-numeric-looking text like "110" or "230" and
# silently convert those cells from text to numbers. Copy/PasteSpecial
# preserves the original cell type exactly, which is what the target
# workbook expects (the group-code/category-name text cells stay text).

$ws.Range("E1:E$lastRow").Copy()
$ws.Range("Z1:Z$lastRow").PasteSpecial()

$ws.Range("F1:F$lastRow").Copy()
$ws.Range("E1:E$lastRow").PasteSpecial()

$ws.Range("Z1:Z$lastRow").Copy()
$ws.Range("F1:F$lastRow").PasteSpecial()

$ws.Range("Z1:Z$lastRow").Clear()
